$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Make room for a new weekly record by shifting every existing data
# row down by one. Inserting at the very bottom (after the last used
# row) lets the new row inherit formatting from the row above it
# (a plain data row) instead of the bold header row, so we copy the
# previous values down from the bottom up rather than inserting at
# row 2 directly.
$ws.Rows("9:9").Insert()
for ($r = 8; $r -ge 2; $r--) {
  for ($c = 1; $c -le 20; $c++) {
    $ws.Cells.Item($r + 1, $c).Value2 = $ws.Cells.Item($r, $c).Value2
  }
}

# Write this week's new record into the now-vacant row 2.
$ws.Range("A2").Value = 7
$ws.Range("B2").Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Range("C2").Value = "Ñuble"
$ws.Range("D2").Value = 44537
$ws.Range("E2").Value = 16
$ws.Range("F2").Value = "Fruta"
$ws.Range("G2").Value = 100103
$ws.Range("H2").Value = "Frutos de hueso (carozo)"
$ws.Range("I2").Value = 100103003
$ws.Range("J2").Value = "Damasco"
$ws.Range("K2").Value = "Castle Brite"
$ws.Range("L2").Value = "Primera"
$ws.Range("M2").Value = 60
$ws.Range("N2").Value = 21000
$ws.Range("O2").Value = 21500
$ws.Range("P2").Value = 21250
$ws.Range("Q2").Value = "`$/caja 15 kilos"
$ws.Range("R2").Value = "Región de O'Higgins"
$ws.Range("S2").Value = 1417
$ws.Range("T2").Value = 15
